$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.315.65'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '1.665.06'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  +0.80%  '
$ws.Range('E5').Value = '  +0.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5349'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.79%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2661'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.65%  '
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.70'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07825'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.566'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').Value = '1.644.85'
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').Value = '1.893.09'
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '0.0₅8245'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.010'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('E19').Value = '  +2.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.82'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('E21').Value = '  +2.07%  '
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('E23').Value = '  +0.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '146.33'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1231'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.201'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.16'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.485'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05845'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.283'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.618'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.282'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.620'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9691'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.824'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.420'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5820'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01604'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8701'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.869'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.70%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '105.29'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.94%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.052.52'
$ws.Range('E42').Value = '  +2.69%  '
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('D44').Value = '1.804.33'
$ws.Range('E44').Value = '  +0.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.88'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.47%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈107'
$ws.Range('E46').Value = '  -4.60%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.015'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.39%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4386'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.051'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05166'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.414'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.55%  '
